$wb = $excel.ActiveWorkbook

# --- Part 1: new "2022-Q1" fund-holdings sheet, inserted right before "总计" ---
$tpl = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"
$tpl.Range("A1:H5").Copy($q1.Range("A1"))
$q1.Range("A1").ClearContents()
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1

$q1.Range("B2").Value = "'012348"
$q1.Range("B2").Style = "Normal"
$q1.Range("C2").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）A"
$q1.Range("D2").Value = "'38.10"
$q1.Range("D2").Style = "Normal"
$q1.Range("E2").Value = "'92.34"
$q1.Range("E2").Style = "Normal"
$q1.Range("F2").Value = "'5.53"
$q1.Range("F2").Style = "Normal"
$q1.Range("G2").Value = "'2.1069"
$q1.Range("G2").Style = "Normal"
$q1.Range("H2").Value = 7

$q1.Range("B3").Value = "'012349"
$q1.Range("B3").Style = "Normal"
$q1.Range("C3").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）C"
$q1.Range("D3").Value = "'14.77"
$q1.Range("D3").Style = "Normal"
$q1.Range("E3").Value = "'92.34"
$q1.Range("E3").Style = "Normal"
$q1.Range("F3").Value = "'5.53"
$q1.Range("F3").Style = "Normal"
$q1.Range("G3").Value = "'0.8168"
$q1.Range("G3").Style = "Normal"
$q1.Range("H3").Value = 7

$q1.Range("B4").Value = "'002379"
$q1.Range("B4").Style = "Normal"
$q1.Range("C4").Value = "工银瑞信香港中小盘股票（QDII）人民币"
$q1.Range("D4").Value = "'1.84"
$q1.Range("D4").Style = "Normal"
$q1.Range("E4").Value = "'86.48"
$q1.Range("E4").Style = "Normal"
$q1.Range("F4").Value = "'4.08"
$q1.Range("F4").Style = "Normal"
$q1.Range("G4").Value = "'0.0751"
$q1.Range("G4").Style = "Normal"
$q1.Range("H4").Value = 5

$q1.Range("B5").Value = "'002380"
$q1.Range("B5").Style = "Normal"
$q1.Range("C5").Value = "工银瑞信香港中小盘股票（QDII）美元"
$q1.Range("D5").Value = "'1.84"
$q1.Range("D5").Style = "Normal"
$q1.Range("E5").Value = "'86.48"
$q1.Range("E5").Style = "Normal"
$q1.Range("F5").Value = "'4.08"
$q1.Range("F5").Style = "Normal"
$q1.Range("G5").Value = "'0.0751"
$q1.Range("G5").Style = "Normal"
$q1.Range("H5").Value = 5

# --- Part 2: prepend a "2022-Q1" row to the "总计" summary sheet ---
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()
$tot.Range("B2:D2").Style = "Normal"
$tot.Range("A3").Copy($tot.Range("A2"))
$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 4
$tot.Range("D2").Value = 3.07
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5

# restore original active sheet/selection
$wb.Worksheets.Item("2020-Q4").Activate()
